# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns for
# rows 2-51 on Sheet1.
#
# The source values are plain text (prices such as "29.433.88" or "1.0000"
# and percentages such as "  -0.50%  "), not real numbers, so when a value
# looks numeric we briefly force a text NumberFormat before writing it and
# then restore the cell to the default "Normal" style afterwards so we don't
# leave a stray style index behind on cells that had none originally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $addr, $value) {
    $range = $sheet.Range($addr)
    if ($value -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}

Set-CellText $ws "D2" '29.433.88'
Set-CellText $ws "E2" '  -0.50%  '
Set-CellText $ws "D3" '1.850.97'
Set-CellText $ws "E3" '  -0.22%  '
Set-CellText $ws "D4" '0.9989'
Set-CellText $ws "E4" '  +0.04%  '
Set-CellText $ws "D5" '241.24'
Set-CellText $ws "D6" '0.6337'
Set-CellText $ws "E6" '  -1.69%  '
Set-CellText $ws "D7" '0.9999'
Set-CellText $ws "E7" '  +0.05%  '
Set-CellText $ws "D8" '4.485.57'
Set-CellText $ws "E8" '  +134.70%  '
Set-CellText $ws "D9" '4.579.58'
Set-CellText $ws "E9" '  +111.07%  '
Set-CellText $ws "D10" '0.07575'
Set-CellText $ws "E10" '  +0.65%  '
Set-CellText $ws "D11" '0.2964'
Set-CellText $ws "E11" '  -1.78%  '
Set-CellText $ws "E12" '  +0.78%  '
Set-CellText $ws "D13" '0.07729'
Set-CellText $ws "E13" '  +0.78%  '
Set-CellText $ws "E14" '  -1.10%  '
Set-CellText $ws "E15" '  -0.86%  '
Set-CellText $ws "E16" '  -1.17%  '
Set-CellText $ws "D17" '0.000009910'
Set-CellText $ws "E17" '  +3.41%  '
Set-CellText $ws "D18" '6.182'
Set-CellText $ws "E18" '  -0.81%  '
Set-CellText $ws "D19" '29.459.98'
Set-CellText $ws "E19" '  -0.49%  '
Set-CellText $ws "D20" '232.02'
Set-CellText $ws "E20" '  -2.44%  '
Set-CellText $ws "E21" '  -1.10%  '
Set-CellText $ws "E22" '  +0.01%  '
Set-CellText $ws "D23" '7.603'
Set-CellText $ws "E23" '  -1.76%  '
Set-CellText $ws "D24" '1.0000'
Set-CellText $ws "E24" '  +0.02%  '
Set-CellText $ws "D25" '155.83'
Set-CellText $ws "E25" '  -1.00%  '
Set-CellText $ws "D26" '0.1387'
Set-CellText $ws "E26" '  -2.02%  '
Set-CellText $ws "D27" '8.401'
Set-CellText $ws "E27" '  -1.53%  '
Set-CellText $ws "E28" '  -0.86%  '
Set-CellText $ws "D29" '4.644.52'
Set-CellText $ws "E29" '  +123.58%  '
Set-CellText $ws "D30" '1.469'
Set-CellText $ws "E30" '  -1.69%  '
Set-CellText $ws "D31" '0.05759'
Set-CellText $ws "E31" '  -3.64%  '
Set-CellText $ws "E32" '  +0.15%  '
Set-CellText $ws "E33" '  -0.23%  '
Set-CellText $ws "D34" '4.018'
Set-CellText $ws "E34" '  -1.63%  '
Set-CellText $ws "D35" '1.856'
Set-CellText $ws "E35" '  -1.53%  '
Set-CellText $ws "E36" '  -1.52%  '
Set-CellText $ws "D37" '0.7170'
Set-CellText $ws "E37" '  -0.92%  '
Set-CellText $ws "E38" '  -0.16%  '
Set-CellText $ws "D39" '1.252.08'
Set-CellText $ws "E39" '  +3.51%  '
Set-CellText $ws "D40" '2.803'
Set-CellText $ws "E40" '  +0.58%  '
Set-CellText $ws "D41" '0.01806'
Set-CellText $ws "E41" '  +1.56%  '
Set-CellText $ws "D42" '0.9041'
Set-CellText $ws "E42" '  -1.17%  '
Set-CellText $ws "D43" '6.109'
Set-CellText $ws "E43" '  -1.39%  '
Set-CellText $ws "E44" '  +0.01%  '
Set-CellText $ws "D45" '101.71'
Set-CellText $ws "E45" '  -0.27%  '
Set-CellText $ws "E46" '  -0.47%  '
Set-CellText $ws "D47" '7.151'
Set-CellText $ws "E47" '  -2.99%  '
Set-CellText $ws "D48" '9.214'
Set-CellText $ws "E48" '  +0.58%  '
Set-CellText $ws "D49" '0.4023'
Set-CellText $ws "E49" '  -1.06%  '
Set-CellText $ws "D50" '1.686'
Set-CellText $ws "E50" '  +1.17%  '
Set-CellText $ws "E51" '  -0.43%  '
